# Scheduled market-data refresh: updates the currentAveragePrice* / LevePrice* /
# LeveProfit* columns (H:N) on several leve rows across the ALC/ARM/CRP/CUL/GSM/LTW/WVR
# sheets, per the latest pull from the Universalis price feed.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 203.75
$ws.Range("I9").Value = 213.33333
$ws.Range("J9").Value = 175
$ws.Range("K9").Value = 213.33333
$ws.Range("L9").Value = 175
$ws.Range("M9").Value = -44.33332999999999
$ws.Range("N9").Value = -513

$ws.Range("H19").Value = 397.15625
$ws.Range("J19").Value = 339.5
$ws.Range("L19").Value = 339.5
$ws.Range("N19").Value = -689.5

$ws.Range("H40").Value = 1711.6666
$ws.Range("I40").Value = 2399.6667
$ws.Range("K40").Value = 2399.6667
$ws.Range("M40").Value = -2224.6667

$ws.Range("H96").Value = 481.9
$ws.Range("J96").Value = 942.6667
$ws.Range("L96").Value = 2828.0001
$ws.Range("N96").Value = -5574.0001

$ws.Range("H100").Value = 1782.1428
$ws.Range("I100").Value = 1579.8334
$ws.Range("K100").Value = 1579.8334
$ws.Range("M100").Value = -1038.8334

$ws.Range("H138").Value = 2851.6667
$ws.Range("I138").Value = 2942.889
$ws.Range("J138").Value = 2841.8928
$ws.Range("K138").Value = 8828.667000000001
$ws.Range("L138").Value = 8525.678400000001
$ws.Range("M138").Value = -3688.667000000001
$ws.Range("N138").Value = -18805.6784

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15785.639
$ws.Range("I32").Value = 10629.08
$ws.Range("J32").Value = 23598.605
$ws.Range("K32").Value = 10629.08
$ws.Range("L32").Value = 23598.605
$ws.Range("M32").Value = -10342.08
$ws.Range("N32").Value = -24172.605

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 2101.8
$ws.Range("J4").Value = 2101.8
$ws.Range("L4").Value = 2101.8
$ws.Range("N4").Value = -2325.8

$ws.Range("H16").Value = 62501280
$ws.Range("I16").Value = 71429820
$ws.Range("J16").Value = 1500
$ws.Range("K16").Value = 71429820
$ws.Range("L16").Value = 1500
$ws.Range("M16").Value = -71429533
$ws.Range("N16").Value = -2074

$ws.Range("H31").Value = 1762.9038
$ws.Range("I31").Value = 1720.6809
$ws.Range("J31").Value = 2159.8
$ws.Range("K31").Value = 1720.6809
$ws.Range("L31").Value = 2159.8
$ws.Range("M31").Value = -1425.6809
$ws.Range("N31").Value = -2749.8

$ws.Range("H34").Value = 1762.9038
$ws.Range("I34").Value = 1720.6809
$ws.Range("J34").Value = 2159.8
$ws.Range("K34").Value = 1720.6809
$ws.Range("L34").Value = 2159.8
$ws.Range("M34").Value = -1518.6809
$ws.Range("N34").Value = -2563.8

$ws.Range("H105").Value = 834.3077
$ws.Range("I105").Value = 745
$ws.Range("J105").Value = 977.2
$ws.Range("K105").Value = 745
$ws.Range("L105").Value = 977.2
$ws.Range("M105").Value = 1002
$ws.Range("N105").Value = -4471.2

$ws.Range("H113").Value = 62501280
$ws.Range("I113").Value = 71429820
$ws.Range("J113").Value = 1500
$ws.Range("K113").Value = 71429820
$ws.Range("L113").Value = 1500
$ws.Range("M113").Value = -71427650
$ws.Range("N113").Value = -5840

$ws.Range("H135").Value = 38701.25
$ws.Range("J135").Value = 38701.25
$ws.Range("L135").Value = 38701.25
$ws.Range("N135").Value = -48841.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1557
$ws.Range("I5").Value = 1557
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 4671
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -4559
$ws.Range("N5").ClearContents()

$ws.Range("H23").Value = 521.5
$ws.Range("J23").Value = 521.5
$ws.Range("L23").Value = 1564.5
$ws.Range("N23").Value = -2034.5

$ws.Range("H32").Value = 2850
$ws.Range("I32").Value = 3200
$ws.Range("J32").Value = 2500
$ws.Range("K32").Value = 9600
$ws.Range("L32").Value = 7500
$ws.Range("M32").Value = -9317
$ws.Range("N32").Value = -8066

$ws.Range("H122").Value = 1029.7291
$ws.Range("I122").Value = 803.1667
$ws.Range("J122").Value = 1105.25
$ws.Range("K122").Value = 7228.5003
$ws.Range("L122").Value = 9947.25
$ws.Range("M122").Value = -4778.5003
$ws.Range("N122").Value = -14847.25

$ws.Range("H131").Value = 34535480
$ws.Range("J131").Value = 72649.234
$ws.Range("L131").Value = 217947.702
$ws.Range("N131").Value = -228027.702

$ws.Range("H135").Value = 1557
$ws.Range("I135").Value = 1557
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 14013
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -11478
$ws.Range("N135").ClearContents()

$ws.Range("H140").Value = 25360.727
$ws.Range("I140").Value = 52371.85
$ws.Range("K140").Value = 157115.55
$ws.Range("M140").Value = -151935.55

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2619.394
$ws.Range("I102").Value = 1653.8636
$ws.Range("K102").Value = 1653.8636
$ws.Range("M102").Value = -31.86359999999991

$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 865.35486
$ws.Range("I16").Value = 867.1923
$ws.Range("J16").Value = 855.8
$ws.Range("K16").Value = 867.1923
$ws.Range("L16").Value = 855.8
$ws.Range("M16").Value = -697.1923
$ws.Range("N16").Value = -1195.8

$ws.Range("H46").Value = 1863
$ws.Range("J46").Value = 1863
$ws.Range("L46").Value = 1863
$ws.Range("N46").Value = -2239

$ws.Range("H61").Value = 1871.2858
$ws.Range("I61").Value = 1439.8
$ws.Range("K61").Value = 1439.8
$ws.Range("M61").Value = -1237.8

$ws.Range("H93").Value = 945.55554
$ws.Range("I93").Value = 941.4
$ws.Range("J93").Value = 966.3333
$ws.Range("K93").Value = 941.4
$ws.Range("L93").Value = 966.3333
$ws.Range("M93").Value = 306.6
$ws.Range("N93").Value = -3462.3333

$ws.Range("H110").Value = 35000
$ws.Range("J110").Value = 35000
$ws.Range("L110").Value = 35000
$ws.Range("N110").Value = -43180

$ws.Range("H113").Value = 1871.2858
$ws.Range("I113").Value = 1439.8
$ws.Range("K113").Value = 1439.8
$ws.Range("M113").Value = 730.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 667.625
$ws.Range("I113").Value = 363.9091
$ws.Range("J113").Value = 1335.8
$ws.Range("K113").Value = 1091.7273
$ws.Range("L113").Value = 4007.4
$ws.Range("M113").Value = 1078.2727
$ws.Range("N113").Value = -8347.4

$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()

$ws.Range("H122").Value = 22729856
$ws.Range("I122").Value = 27780370
$ws.Range("J122").Value = 2550
$ws.Range("K122").Value = 83341110
$ws.Range("L122").Value = 7650
$ws.Range("M122").Value = -83338660
$ws.Range("N122").Value = -12550

$ws.Range("H132").Value = 11143.556
$ws.Range("I132").Value = 11143.556
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 33430.66800000001
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -30900.66800000001
$ws.Range("N132").ClearContents()

$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0

$ws.Range("H136").Value = 1763.6364
$ws.Range("I136").Value = 1300
$ws.Range("J136").Value = 2150
$ws.Range("K136").Value = 3900
$ws.Range("L136").Value = 6450
$ws.Range("M136").Value = -1350
$ws.Range("N136").Value = -11550

$ws.Range("H141").Value = 77413.336
$ws.Range("J141").Value = 77413.336
$ws.Range("L141").Value = 77413.336
$ws.Range("N141").Value = -87773.336
